# pas533_newVinAdded small cleaning NewVIN_UT_SS.xlsx : VIN version changed to SYMBOL_2017
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# VERSION column (B2:B5): SYMBOL_2000 -> SYMBOL_2017
$ws.Range("B2:B5").Value = "SYMBOL_2017"

# BI_SYMBOL / PD_SYMBOL / UM_SYMBOL / MP_SYMBOL columns (AC:AF) for rows 2-5
# Row 2 (previously all "C")
$ws.Range("AC2").Value = "BI001"
$ws.Range("AD2").Value = "PD001"
$ws.Range("AE2").Value = "UM001"
$ws.Range("AF2").Value = "MP001"

# Row 3 (previously all "N")
$ws.Range("AC3").Value = "BI002"
$ws.Range("AD3").Value = "PD002"
$ws.Range("AE3").Value = "UM002"
$ws.Range("AF3").Value = "MP002"

# Row 4 (previously all "K")
$ws.Range("AC4").Value = "BI003"
$ws.Range("AD4").Value = "PD003"
$ws.Range("AE4").Value = "UM003"
$ws.Range("AF4").Value = "MP003"

# Row 5 (previously all "N")
$ws.Range("AC5").Value = "BI004"
$ws.Range("AD5").Value = "PD004"
$ws.Range("AE5").Value = "UM004"
$ws.Range("AF5").Value = "MP004"

# Update the view: drop the frozen/scrolled topLeftCell="V1" and move selection to E9
$ws.Range("E9").Select()
